$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.049101265588916
$ws.Range("D2").Value = 0.1924808555325725
$ws.Range("E2").Value = 0.2277992248769998
$ws.Range("F2").Value = 1.882514647723873
$ws.Range("G2").Value = 0.002498220423063922
$ws.Range("J2").Value = 0.3146177846317011
$ws.Range("K2").Value = 0.4012938604552971
$ws.Range("L2").Value = 0.1503062103076473
$ws.Range("O2").Value = 4.671985175936385
$ws.Range("B3").Value = 1.023879288246235
$ws.Range("D3").Value = 0.1922032250418582
$ws.Range("E3").Value = 0.2292256054191153
$ws.Range("F3").Value = 1.891121815335779
$ws.Range("G3").Value = 0.002500896661608302
$ws.Range("J3").Value = 0.3170678202612676
$ws.Range("K3").Value = 0.3523913072544929
$ws.Range("L3").Value = 0.1404689338976084
$ws.Range("O3").Value = 4.697265035710018
$ws.Range("B4").Value = 1.008803515631655
$ws.Range("D4").Value = 0.1920927122406439
$ws.Range("E4").Value = 0.2301681925408392
$ws.Range("F4").Value = 1.897281754352413
$ws.Range("G4").Value = 0.002502628719015291
$ws.Range("J4").Value = 0.3186582457664624
$ws.Range("K4").Value = 0.3223164720430418
$ws.Range("L4").Value = 0.1344635467498136
$ws.Range("O4").Value = 4.715064530342488
$ws.Range("B5").Value = 1.002763880005034
$ws.Range("D5").Value = 0.1920628071258719
$ws.Range("E5").Value = 0.2305691276605781
$ws.Range("F5").Value = 1.900012186019495
$ws.Range("G5").Value = 0.002503356951225777
$ws.Range("J5").Value = 0.3193280369804636
$ws.Range("K5").Value = 0.3100492558761516
$ws.Range("L5").Value = 0.1320252083475424
$ws.Range("O5").Value = 4.722890870556171
$ws.Range("B6").Value = 1.001767294220059
$ws.Range("D6").Value = 0.1920587567926901
$ws.Range("E6").Value = 0.2306367195916081
$ws.Range("F6").Value = 1.90047887659874
$ws.Range("G6").Value = 0.002503479228763695
$ws.Range("J6").Value = 0.3194405656245936
$ws.Range("K6").Value = 0.3080116218447984
$ws.Range("L6").Value = 0.1316208670070438
$ws.Range("O6").Value = 4.724225035871541
$ws.Range("B7").Value = 1.008721641641529
$ws.Range("D7").Value = 0.1920922475974507
$ws.Range("E7").Value = 0.2301735315355415
$ws.Range("F7").Value = 1.897317686097615
$ws.Range("G7").Value = 0.002502638449381906
$ws.Range("J7").Value = 0.3186671909824002
$ws.Range("K7").Value = 0.3221510775974821
$ws.Range("L7").Value = 0.134430626180432
$ws.Range("O7").Value = 4.715167759200028
$ws.Range("B8").Value = 1.040319858066056
$ws.Range("D8").Value = 0.1923727171277321
$ws.Range("E8").Value = 0.2282771979270919
$ws.Range("F8").Value = 1.885300862294073
$ws.Range("G8").Value = 0.002499124792152703
$ws.Range("J8").Value = 0.3154447069543256
$ws.Range("K8").Value = 0.3844427308827107
$ws.Range("L8").Value = 0.1469072055153333
$ws.Range("O8").Value = 4.680229074780272
$ws.Range("B9").Value = 1.10551927782825
$ws.Range("D9").Value = 0.1933962227372987
$ws.Range("E9").Value = 0.225087059275296
$ws.Range("F9").Value = 1.868673689125885
$ws.Range("G9").Value = 0.002492936371468923
$ws.Range("J9").Value = 0.3098070876127057
$ws.Range("K9").Value = 0.5061866232550472
$ws.Range("L9").Value = 0.1716432306670015
$ws.Range("O9").Value = 4.629779819058029
$ws.Range("B10").Value = 1.155369187544437
$ws.Range("D10").Value = 0.1944341541984969
$ws.Range("E10").Value = 0.2230636722501611
$ws.Range("F10").Value = 1.860680927304443
$ws.Range("G10").Value = 0.002488813370109475
$ws.Range("J10").Value = 0.3060785961477785
$ws.Range("K10").Value = 0.5953568830431664
$ws.Range("L10").Value = 0.1899747243146521
$ws.Range("O10").Value = 4.603723378963593
$ws.Range("B11").Value = 1.178465370462987
$ws.Range("D11").Value = 0.1949678758126296
$ws.Range("E11").Value = 0.2222123696389993
$ws.Range("F11").Value = 1.857960639214639
$ws.Range("G11").Value = 0.00248702879671708
$ws.Range("J11").Value = 0.3044717425038312
$ws.Range("K11").Value = 0.6358582956055727
$ws.Range("L11").Value = 0.1983472701045343
$ws.Range("O11").Value = 4.594259377394224
$ws.Range("B12").Value = 1.187271047058402
$ws.Range("D12").Value = 0.195178784077207
$ws.Range("E12").Value = 0.2218999166828315
$ws.Range("F12").Value = 1.857062085888728
$ws.Range("G12").Value = 0.002486366043117957
$ws.Range("J12").Value = 0.303876071918503
$ws.Range("K12").Value = 0.6511855355721252
$ws.Range("L12").Value = 0.2015224025355025
$ws.Range("O12").Value = 4.591019070031138
$ws.Range("B13").Value = 1.185371943937156
$ws.Range("D13").Value = 0.1951329706919154
$ws.Range("E13").Value = 0.2219667683797404
$ws.Range("F13").Value = 1.8572497559813
$ws.Range("G13").Value = 0.002486508200648881
$ws.Range("J13").Value = 0.3040037909019375
$ws.Range("K13").Value = 0.6478849867071688
$ws.Range("L13").Value = 0.2008383778482283
$ws.Range("O13").Value = 4.591701651702493
$ws.Range("B14").Value = 1.17918862682663
$ws.Range("D14").Value = 0.1949850513104465
$ws.Range("E14").Value = 0.2221864653475798
$ws.Range("F14").Value = 1.857884078650173
$ws.Range("G14").Value = 0.00248697401101119
$ws.Range("J14").Value = 0.3044224797363428
$ws.Range("K14").Value = 0.6371194769611179
$ws.Range("L14").Value = 0.1986083983764075
$ws.Range("O14").Value = 4.593985911360249
$ws.Range("B15").Value = 1.175408916568927
$ws.Range("D15").Value = 0.1948955908002432
$ws.Range("E15").Value = 0.2223223267744867
$ws.Range("F15").Value = 1.858289749420479
$ws.Range("G15").Value = 0.002487261027649936
$ws.Range("J15").Value = 0.3046806062448182
$ws.Range("K15").Value = 0.6305239969043441
$ws.Range("L15").Value = 0.1972430693280671
$ws.Range("O15").Value = 4.595429819002277
$ws.Range("B16").Value = 1.153868179474529
$ws.Range("D16").Value = 0.1944005084755318
$ws.Range("E16").Value = 0.2231206960288983
$ws.Range("F16").Value = 1.860877120545766
$ws.Range("G16").Value = 0.002488931822081308
$ws.Range("J16").Value = 0.3061854017374168
$ws.Range("K16").Value = 0.5927086955047116
$ws.Range("L16").Value = 0.1894282165934129
$ws.Range("O16").Value = 4.604389945509809
$ws.Range("B17").Value = 1.140760556780435
$ws.Range("D17").Value = 0.1941125211072858
$ws.Range("E17").Value = 0.2236281609976789
$ws.Range("F17").Value = 1.862698833153729
$ws.Range("G17").Value = 0.002489980064645471
$ws.Range("J17").Value = 0.3071313878290791
$ws.Range("K17").Value = 0.569493656745351
$ws.Range("L17").Value = 0.1846424984210984
$ws.Range("O17").Value = 4.610498582987162
$ws.Range("B18").Value = 1.133260894270251
$ws.Range("D18").Value = 0.1939526759982186
$ws.Range("E18").Value = 0.2239265513685194
$ws.Range("F18").Value = 1.863832826850867
$ws.Range("G18").Value = 0.002490591555139481
$ws.Range("J18").Value = 0.3076838966618594
$ws.Range("K18").Value = 0.5561351351332462
$ws.Range("L18").Value = 0.1818930397903671
$ws.Range("O18").Value = 4.614236996486994
$ws.Range("B19").Value = 1.130728439952264
$ws.Range("D19").Value = 0.1938995525513292
$ws.Range("E19").Value = 0.22402870011231
$ws.Range("F19").Value = 1.864231584673234
$ws.Range("G19").Value = 0.002490800068810766
$ws.Range("J19").Value = 0.3078724105841593
$ws.Range("K19").Value = 0.5516111887168336
$ws.Range("L19").Value = 0.1809626685529082
$ws.Range("O19").Value = 4.615541385406971
$ws.Range("B20").Value = 1.142151802527138
$ws.Range("D20").Value = 0.1941425782083641
$ws.Range("E20").Value = 0.2235734669217955
$ws.Range("F20").Value = 1.862495988899838
$ws.Range("G20").Value = 0.002489867591000033
$ws.Range("J20").Value = 0.307029816532447
$ws.Range("K20").Value = 0.5719655477860783
$ws.Range("L20").Value = 0.1851516205619959
$ws.Range("O20").Value = 4.6098250335844
$ws.Range("B21").Value = 1.181003202970885
$ws.Range("D21").Value = 0.1950282604252251
$ws.Range("E21").Value = 0.2221216660907093
$ws.Range("F21").Value = 1.857694193064006
$ws.Range("G21").Value = 0.002486836838092315
$ws.Range("J21").Value = 0.3042991532131749
$ws.Range("K21").Value = 0.6402818392657537
$ws.Range("L21").Value = 0.1992632731487873
$ws.Range("O21").Value = 4.593305646793254
$ws.Range("B22").Value = 1.206742254273252
$ws.Range("D22").Value = 0.1956583705378492
$ws.Range("E22").Value = 0.221230621283425
$ws.Range("F22").Value = 1.855322727215153
$ws.Range("G22").Value = 0.002484931960525472
$ws.Range("J22").Value = 0.3025891620782994
$ws.Range("K22").Value = 0.6848732086495772
$ws.Range("L22").Value = 0.208512964150259
$ws.Range("O22").Value = 4.584511451907929
$ws.Range("B23").Value = 1.192973256827344
$ws.Range("D23").Value = 0.1953173955717133
$ws.Range("E23").Value = 0.2217009094488738
$ws.Range("F23").Value = 1.85651829871135
$ws.Range("G23").Value = 0.002485941705102353
$ws.Range("J23").Value = 0.3034949937305118
$ws.Range("K23").Value = 0.6610794485155793
$ws.Range("L23").Value = 0.2035738269920984
$ws.Range("O23").Value = 4.589021899773314
$ws.Range("B24").Value = 1.141522707879318
$ws.Range("D24").Value = 0.1941289715645595
$ws.Range("E24").Value = 0.2235981734226726
$ws.Range("F24").Value = 1.862587424815302
$ws.Range("G24").Value = 0.002489918412869226
$ws.Range("J24").Value = 0.3070757099734607
$ws.Range("K24").Value = 0.5708480428164364
$ws.Range("L24").Value = 0.1849214404661694
$ws.Range("O24").Value = 4.610128839796801
$ws.Range("B25").Value = 1.087537274237974
$ws.Range("D25").Value = 0.1930689250258695
$ws.Range("E25").Value = 0.2258936782250114
$ws.Range("F25").Value = 1.8724296864252
$ws.Range("G25").Value = 0.002494535804938814
$ws.Range("J25").Value = 0.3112594535078248
$ws.Range("K25").Value = 0.4732981652443016
$ws.Range("L25").Value = 0.1649232895029655
$ws.Range("O25").Value = 4.641494119915365
